$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.157.33'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '2.500.90'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.92'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.62'
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.526'
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.538'
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.72'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.23'
$ws.Range('E11').Value = '  +8.36%  '
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.14'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = '2.894.52'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').Value = '2.503.86'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.839'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').Value = '48.011.46'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.77'
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.73'
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '278.04'
$ws.Range('E23').Value = '  +12.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.96'
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.69'
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.75'
$ws.Range('E28').Value = '  -2.78%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.140'
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.12'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.11'
$ws.Range('E31').Value = '  -4.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.42'
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.70'
$ws.Range('E33').Value = '  -2.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.33'
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('E36').Value = '  -1.22%  '
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('E39').Value = '  -1.75%  '
$ws.Range('E40').Value = '  -0.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '121.74'
$ws.Range('E41').Value = '  +1.47%  '
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.47'
$ws.Range('E43').Value = '  -4.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0301'
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('D45').Value = '2.014.71'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.17'
$ws.Range('E46').Value = '  +3.90%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.85'
$ws.Range('E47').Value = '  +2.01%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.00'
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.01'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.29'
$ws.Range('E51').Value = '  +2.83%  '
